# Apply the 2024-03-01 cryptos-list refresh (GitHub Actions data update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing Text storage (column D holds prices
# like "1.00" / "3.444.25" that must stay text, not get reinterpreted as
# numbers). Style is restored afterward so no stray formatting is left.
function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextCell "D2" '62.582.85'
$ws.Range("E2").Value = '  +0.64%  '

Set-TextCell "D3" '3.428.30'
$ws.Range("E3").Value = '  +0.61%  '

$ws.Range("E4").Value = '  +0.17%  '

Set-TextCell "D5" '406.94'
$ws.Range("E5").Value = '  -0.02%  '

Set-TextCell "D6" '130.33'
$ws.Range("E6").Value = '  -0.39%  '

$ws.Range("E7").Value = '  -2.08%  '

$ws.Range("E8").Value = '  +0.04%  '

Set-TextCell "D9" '0.692'
$ws.Range("E9").Value = '  +2.11%  '

Set-TextCell "D10" '0.137'
$ws.Range("E10").Value = '  +7.97%  '

Set-TextCell "D11" '41.94'
$ws.Range("E11").Value = '  -0.67%  '

$ws.Range("E12").Value = '  -0.15%  '

Set-TextCell "D13" '8.39'
$ws.Range("E13").Value = '  -2.20%  '

Set-TextCell "D14" '19.81'
$ws.Range("E14").Value = '  +0.06%  '

Set-TextCell "D15" '3.407.03'
$ws.Range("E15").Value = '  +0.34%  '

Set-TextCell "D16" '62.496.57'
$ws.Range("E16").Value = '  +0.75%  '

Set-TextCell "D17" '11.56'
$ws.Range("E17").Value = '  -0.10%  '

$ws.Range("E18").Value = '  -1.81%  '

Set-TextCell "D19" '0.0000155'
$ws.Range("E19").Value = '  +15.38%  '

$ws.Range("E20").Value = '  -3.29%  '

Set-TextCell "D21" '84.30'
$ws.Range("E21").Value = '  +1.83%  '

Set-TextCell "D22" '312.33'
$ws.Range("E22").Value = '  +1.30%  '

Set-TextCell "D23" '12.78'
$ws.Range("E23").Value = '  -3.19%  '

Set-TextCell "D24" '3.17'
$ws.Range("E24").Value = '  +0.07%  '

Set-TextCell "D25" '4.76'
$ws.Range("E25").Value = '  +1.68%  '

Set-TextCell "D26" '29.69'
$ws.Range("E26").Value = '  -0.34%  '

$ws.Range("E27").Value = '  -6.06%  '

$ws.Range("E28").Value = '  +3.95%  '

$ws.Range("E29").Value = '  +4.38%  '

Set-TextCell "D30" '44.63'
$ws.Range("E30").Value = '  +5.07%  '

$ws.Range("E32").Value = '  -1.59%  '

Set-TextCell "D33" '11.37'
$ws.Range("E33").Value = '  -3.68%  '

$ws.Range("E34").Value = '  +0.12%  '

$ws.Range("E35").Value = '  -1.38%  '

$ws.Range("E36").Value = '  -1.23%  '

Set-TextCell "D37" '1.00'
$ws.Range("E37").Value = '  +0.26%  '

$ws.Range("B38").Value = 'TheGraph'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell "D38" '0.325'
$ws.Range("E38").Value = '  +13.51%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell "D39" '2.97'
$ws.Range("E39").Value = '  -0.46%  '

$ws.Range("E40").Value = '  -4.26%  '

Set-TextCell "D41" '142.19'
$ws.Range("E41").Value = '  +3.35%  '

$ws.Range("E42").Value = '  -0.32%  '

$ws.Range("E43").Value = '  -2.61%  '

Set-TextCell "D44" '3.93'

Set-TextCell "D45" '16.78'
$ws.Range("E45").Value = '  -1.89%  '

$ws.Range("E46").Value = '  -0.63%  '

Set-TextCell "D47" '21.21'
$ws.Range("E47").Value = '  -2.37%  '

Set-TextCell "D48" '2.102.57'
$ws.Range("E48").Value = '  -2.35%  '

$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell "D49" '2.32'
$ws.Range("E49").Value = '  -1.26%  '

$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextCell "D50" '1.96'
$ws.Range("E50").Value = '  +2.47%  '

Set-TextCell "D51" '1.10'
$ws.Range("E51").Value = '  +29.67%  '
